# Commit: "add monthly yearly functionality in customers and tax sections,
#          add migration for add tax_duration col"
#
# The workbook-level change behind that migration is a new "Tax Duration"
# header column appended to the company sheet's header row (column Q,
# right after the existing "Industry" column in P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell -> new shared string "Tax Duration", dimension grows to
# A1:Q1 and row 1's span becomes 1:17 automatically.
$ws.Range("Q1").Value = "Tax Duration"

# Match the author's resulting selection state (active cell moves to O1
# after the edit).
$ws.Range("O1").Select() | Out-Null
